$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @(
    "darryl.fenner@arcanite-ndt.com",
    "h_tazikeh@yahoo.com",
    "ben.li@mistrasgroup.com",
    "jledden@metalogicinspection.com",
    "robert.hindle@brucepower.com",
    "peter.flaman@enbridge.com",
    "darrell.skinner5@gmail.com",
    "huibin.hu@brucepower.com",
    "helderassuncao@hotmail.com",
    "brunobandeiram@hotmail.com",
    "amin.maki@gmail.com",
    "dreimer.cgy@gmail.com",
    "rickertjohn@icloud.com",
    "david199130@gmail.com",
    "dale.chadney@enbridge.com",
    "Ybehnamian@gmail.com",
    "amirbehvandi01747@gmail.com",
    "elohorakiri@gmail.com",
    "lcsouza97@gmail.com",
    "mungdx66@gmail.com",
    "logan.toth@enbridge.com",
    "brook.althouse@arcanite-ndt.com",
    "luay.ahmed89@yahoo.com",
    "josh_spencer@transcanada.com",
    "dstrabel@telus.net",
    "derryckm@gmail.com",
    "carlb.soares@live.com",
    "dheuston@ndtgroup.ca",
    "jhetu@ndtgroup.ca",
    "maurolcardoso@gmail.com",
    "Sohail.anwar@cnrl.com",
    "Gengsheng.Weng@cnrl.com",
    "twillier@metalogicinspection.com",
    "roger.fourny@shaw.ca",
    "k.ndt2008@gmail.com",
    "mfc_98@hotmail.com",
    "mwgumb@bwxt.com",
    "justin.knutsen@metalogicinspection.com",
    "gford@metalogicinspection.com",
    "lebarrozo@gmail.com",
    "janalizadeh@gmail.com",
    "Nathan.Schuler@cnrl.com",
    "elaflamme@nucleom.ca",
    "naldenir.amaral@gmail.com",
    "dalsaunders@gmail.com",
    "natanoel04@hotmail.com",
    "marco.venne@gmail.com",
    "ideraldo.tiburcio@sbmoffshore.com",
    "james.nimijohn@enbridge.com",
    "Alex.Arrau@cnrl.com",
    "Mandy.Nelson@cnrl.com",
    "lucioinspetordm2@gmail.com",
    "johnyboyaz@gmail.com",
    "marcelo.borchert@sbmoffshore.com",
    "robertopintocq@gmail.com",
    "Tylor.Arguin@wav.ca",
    "Mike.Cook@sbdinc.com",
    "Ryan.Faubert@enbridge.com",
    "Udaya.Sundar@cnrl.com",
    "matthew.prowse@acuren.com",
    "homayoun.javadi@tescan.ca",
    "David-Tompkins@hotmail.com",
    "Mojtaba.ghaderi60@gmail.com",
    "noelson.amaral@gmail.com",
    "timnelson155@msn.com",
    "Michael.Brault@irisndt.com",
    "VPopov@nucleom.ca",
    "ssinger@acuren.com",
    "belchiorvirgilio@gmail.com",
    "jodland@tiltinspection.com",
    "renatokow@hotmail.com",
    "marko.alekszity@gmail.com",
    "weston.ellis@opg.com",
    "john.reardon@opg.com",
    "Brian.Purves@wav.ca",
    "kelly.norman@wav.ca",
    "mohammad.koochak@gmail.com",
    "doubrumm@gmail.com",
    "rajasengodan55@gmail.com",
    "rgarcia@ndtgroup.ca",
    "saramella@gmail.com",
    "paul.spencer@applusrtd.com",
    "rsb407@gmail.com",
    "Philippe.Cyr@acuren.com",
    "dana.martin@opg.com",
    "Sheryl.vanderfluit@mistrasgroup.com",
    "amirghabraee@yahoo.com",
    "Jnsmith@bwxt.com",
    "wcpamer@bwxt.com",
    "practitionerinspection@gmail.com",
    "Craig.McMeeken@wav.ca",
    "Ben.Ren@cnrl.com",
    "Muhammad.Akbar@cnrl.com",
    "Hisham.Madi@cnrl.com",
    "jtreacy@its-ndt.com",
    "dinghaifeng8888@hotmail.com",
    "garret.elkins@gmail.com",
    "insptonelli@hotmail.com",
    "nunomcmarques@hotmail.com",
    "Jaeger.Lonsdale@sbdinc.com",
    "levikitt@live.ca",
    "Bradley.Kuntz@enbridge.com",
    "Jay.Brooks@enbridge.com",
    "doug.desruisseaux@enbridge.com",
    "logan.campbell@enbridge.com",
    "darryl.czajkowski@opg.com",
    "j.lawson@opg.com",
    "derrick.watson@brucepower.com",
    "chabouni.djr@gmail.com",
    "owen.nicol@enbridge.com",
    "dale.berezan@applusrtd.com",
    "jian-zhao@hotmail.com",
    "scott.robinson@enbridge.com",
    "b_ilkuf@hotmail.com",
    "douglamarre@sympatico.ca",
    "r.hoffmann@rae.com",
    "Waleed.Rafiq@cnrl.com",
    "Perry.Lawless@cnrl.com",
    "Stephen.Orser@cnrl.com",
    "Qaiser.Butt@cnrl.com",
    "Ayo.Salaudeen@cnrl.com",
    "jonathan.Uhlman@cnrl.com",
    "Anand.Palani@cnrl.com",
    "scott.bangs@brucepower.com",
    "tyler.rickard@enbridge.com",
    "Zawar.Muhammad@cnrl.com",
    "daniel.norman@opg.com",
    "hossein.taheri8067@gmail.com",
    "matt.fritz@arcanite-ndt.com",
    "awallace441@gmail.com",
    "joel.djordjevic@arcanite-ndt.com",
    "cchartier@ndtgroup.ca",
    "humbraganca@hotmail.com",
    "flokinn_coq@hotmail.com",
    "jmay@ndtgroup.ca",
    "marcus.accon@gmail.com",
    "ssusac@ndtgroup.ca",
    "rackitndt@outlook.com",
    "nic.shoebridge@enbridge.com",
    "jsaint@metalogicinspection.com",
    "blake.macpherson@enbridge.com",
    "eng.ndt@yahoo.com",
    "h.herrera@rae.com",
    "Hassan.Sattar@cnrl.com",
    "Marat.Kireev@cnrl.com",
    "Alexie.Broddy@cnrl.com",
    "nima_vakil@yahoo.com",
    "nrweston@anodendt.ca",
    "barry.giasson@opg.com",
    "michelle.fry@opg.com",
    "devon.algera@opg.com",
    "briancable40@gmail.com",
    "jeehmorais89@gmail.com",
    "colbyritzut@gmail.com",
    "Dean.Ikert@enbridge.com",
    "gabrielvfortes@gmail.com",
    "nikhilsatheesan@outlook.com",
    "jesse.rempel@enbridge.com",
    "karandeepgill760@gmail.com",
    "ewhite@ndtgroup.ca",
    "collin.coffey@mantech.com",
    "Curtis.behnke@enbridge.com",
    "marcel@strauhs.com.br",
    "ryan.ziefflie@enbridge.com",
    "carl.gerbrandt@wav.ca",
    "andersonreisssilva@gmail.com",
    "gabriell_008@hotmail.com",
    "nima.vakilotojjar@gmail.com",
    "vinicius.bogos@gmail.com",
    "paul.burton@opg.com",
    "limogesjon@gmail.com",
    "fernando.grigolato@outlook.com",
    "nickolas.lau@brucepower.com",
    "taylor.gardiner@applusrtd.com",
    "srsouzapinto@gmail.com",
    "qualidade.daniel@gmail.com",
    "ferrerinsp@gmail.com",
    "nathaliasg.fr@hotmail.com",
    "Jarratt.Bilodeau@StuartOlson.com",
    "thiago_claro@hotmail.com",
    "julionunescosta@gmail.com",
    "chaboki.ali@gmail.com",
    "a.beckman94@hotmail.com",
    "lacianelli@hotmail.com",
    "eng.douglaswilson@yahoo.com.br",
    "shawn.hanrahan@opg.com",
    "sa-jailson@hotmail.com",
    "blessed.agunu@totalenergies.com",
    "jrhilson16@gmail.com",
    "fqlara@gmail.com",
    "mendoncalon@gmail.com",
    "oko.oono@gmail.com",
    "soheyltahan@gmail.com",
    "gmalmeida@isgbrasil.com.br",
    "yingsong.wu@applusrtd.com",
    "vanderlei_nogueira@hotmail.com",
    "felipe@orionsic.com.br",
    "eric.scott999@gmail.com",
    "brunoinsp@hotmail.com",
    "jackbloy@gmail.com",
    "rodrigostohler@gmail.com",
    "Vptamy@gmail.com",
    "diegom.silva@outlook.com",
    "felipecadiente@gmail.com",
    "marcelo.dimensional@gmail.com",
    "wesleyhweber@gmail.com",
    "hamed.faghihi@opg.com",
    "ben.Leblanc@ultratest.ca",
    "adam.watkins@enbridge.com",
    "kejqual@telusplanet.net",
    "xrayedit23@hotmail.com",
    "frank.santana@opg.com",
    "christopher.wood@opg.com",
    "mark.wainman@opg.com",
    "brock.vangaver@enbridge.com",
    "anderson.ferraz@sbmoffshore.com",
    "robertgwif@yahoo.com.br",
    "shaun.mcassey@ultratest.ca",
    "pedro_augusto27@hotmail.com",
    "alexandre_borchert@hotmail.com",
    "sean.villeneuve@opg.com",
    "psinats@elanderinspection.ca",
    "Paulo.Gaviria@enbridge.com",
    "mattklassen@telus.net",
    "tom.kroeker@enbridge.com",
    "alberto.oliveira@sbmoffshore.com",
    "cforero@ndtgroup.ca",
    "bzieger@metalogicinspection.com",
    "maurilio.filho@sbmoffshore.com",
    "chrisdleslie@outlook.com",
    "John.Adigun@cnrl.com",
    "Syed.Raza@cnrl.com",
    "luis.marka@inphaseintegrity.com",
    "marco.proenca@yahoo.com.br",
    "motasilvab@yahoo.com.br",
    "claudio.quiala@vipaqui-angola.com",
    "sergio.targino@hotmail.com",
    "srquality.inspecao@hotmail.com",
    "shitaleny00@hotmail.com",
    "gabrielsalesmaia@hotmail.com",
    "rthomson@bwxt.com",
    "joseph.konopka@mistrasgroup.com",
    "alex.butcher@brucepower.com",
    "Essam.E.Elnahrawy@lngcanada.ca",
    "calvinkboyle@gmail.com",
    "leonc_rj@yahoo.com.br",
    "cesar.valero@repsol.com",
    "noelfer@gmail.com",
    "greg.braham@caveinspection.com",
    "elia@leveliiindt.com",
    "FBLeclerc@nucleom.ca",
    "sardinhagustavo@hotmail.com",
    "annyscordeiro@gmail.com",
    "ivanelsonogueira36@gmail.com",
    "Saeed.Farea@cnrl.com",
    "petronio@petrobras.com.br",
    "fatihy100@gmail.com",
    "henrique.psouza@hotmail.com",
    "jardelmont@gmail.com",
    "bhardwajmehta22@gmail.com",
    "weiming342@gmail.com",
    "brendanmatthews@mail.com",
    "nima_khatib@yahoo.com",
    "nicolo.mattina@opg.com",
    "jeff.gebhart@evrazna.com",
    "omerzamir3@gmail.com",
    "luisfelipesel@gmail.com",
    "TMorrison@summitinspection.ca",
    "curtis.glen@opg.com",
    "ebx.consultoria@gmail.com",
    "tysonosmond@msn.com",
    "patrick.lannigan@brucepower.com",
    "gina.mc@jgmservices.ca",
    "kellen.daly@opg.com",
    "diego.gomes@me.com",
    "anguxlam@gmail.com",
    "brian.millejours@opg.com"
)

for ($i = 0; $i -lt $emails.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $emails[$i]
}

Write-Host ("wrote " + $emails.Length + " emails")
